$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RESOURCES")

# Insert a new row above current row 2 ("none"/"NONE" resource), shifting
# all the other resource rows down by one.
$ws.Rows.Item(2).Insert()

# Fill in the new row's values.
$ws.Range("A2").Value = "none"
$ws.Range("B2").Value = "NONE"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# Match the formatting used by the rest of the table: column A uses the
# "left aligned" style, columns B:F use the "centered" style.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("B2:F2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("F14").Select()
